$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format
# before assignment so Excel keeps them as strings (matching the
# original inline-string cell type) instead of converting to numbers.
$textCells = @('D5', 'D7', 'D8', 'D9', 'D11', 'D12', 'D13', 'D15', 'D17', 'D18', 'D19', 'D22', 'D24', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D39', 'D43', 'D44', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume values
$ws.Range('D2').Value = '28.520.88'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.820.17'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D5').Value = '317.21'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.5157'
$ws.Range('E7').Value = '  -3.48%  '
$ws.Range('D8').Value = '0.3879'
$ws.Range('E8').Value = '  -2.73%  '
$ws.Range('D9').Value = '0.08421'
$ws.Range('E9').Value = '  +8.26%  '
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').Value = '1.110'
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').Value = '6.434'
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('D13').Value = '20.99'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').Value = '7.514'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '1.822.86'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').Value = '92.80'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '0.06683'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D22').Value = '6.083'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '28.564.98'
$ws.Range('D24').Value = '11.39'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').Value = '21.08'
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').Value = '159.23'
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('D28').Value = '2.031.37'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').Value = '2.412'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '126.00'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '0.1087'
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').Value = '1.094'
$ws.Range('E32').Value = '  -5.21%  '
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').Value = '0.07512'
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('D35').Value = '3.683'
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = '8.742'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  -1.45%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '1.403'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('D44').Value = '13.51'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '3.775'
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '125.98'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '1.990'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').Value = '1.199'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').Value = '0.06976'
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('D51').Value = '74.41'
$ws.Range('E51').Value = '  -0.33%  '
